$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 72, shifting existing rows 72-86 down to 73-87
$ws.Rows.Item(72).Insert()

# Fill in the new row 72 with the new record's data
$ws.Cells.Item(72, 1).Value2 = 1
$ws.Cells.Item(72, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(72, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(72, 4).Value2 = 44722
$ws.Cells.Item(72, 5).Value2 = 15
$ws.Cells.Item(72, 6).Value = "Fruta"
$ws.Cells.Item(72, 7).Value2 = 100102
$ws.Cells.Item(72, 8).Value = "Cítricos"
$ws.Cells.Item(72, 9).Value2 = 100102005
$ws.Cells.Item(72, 10).Value = "Naranja"
$ws.Cells.Item(72, 11).Value = "Fukumoto"
$ws.Cells.Item(72, 12).Value = "Tercera"
$ws.Cells.Item(72, 13).Value2 = 270
$ws.Cells.Item(72, 14).Value2 = 700
$ws.Cells.Item(72, 15).Value2 = 750
$ws.Cells.Item(72, 16).Value2 = 725
$ws.Cells.Item(72, 17).Value = "$/caja 18 kilos importada"
$ws.Cells.Item(72, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(72, 19).Value2 = 40
$ws.Cells.Item(72, 20).Value2 = 18
